$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Картон): update "от 200 до 500кг" and "от 500кг до 1 тонны" prices
$ws.Range("D4").Value = 5.5
$ws.Range("E4").Value = 6

# Row 7 (Маккулатура книги журналы**): update prices and change Самоподвоз
# "свыше 1 тонны" value from a price to "дог", and the address/range column
# (G) from "5-6" to "6-7"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5.5
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = "дог"
$ws.Range("G7").Value = "6-7"

# Row 8 (Газеты): same kind of update as row 7
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5.5
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = "дог"
$ws.Range("G8").Value = "6-7"

# Update the active selection to match the author's final cursor position
$null = $ws.Range("E7").Select()
